# Fixing some formatting typos in 'Dynamic Memory in C.pptx'.
#
# Slide 10 (the "Example function" code listing) has two lines whose
# leading indentation was only 2 spaces while the surrounding code uses a
# 4-space indent:
#
#   "  plotHistogram(histogram, numNumbers);"   ->  "    plotHistogram(histogram, numNumbers);"
#   "  free(histogram); histogram = NULL;"       ->  "    free(histogram); histogram = NULL;"
#
# Both lines live as single runs inside the body placeholder (shape 2) of
# slide 10. We use TextRange.Find to locate each line by its exact current
# text and then overwrite just that sub-range's .Text, which preserves the
# run's existing character formatting (rPr: Courier New latin/ea/cs/sym)
# untouched - only the literal text content changes, matching the diff.
#
# NOTE: Find() always searches from the start of the TextRange and returns
# a freshly-located TextRange, so doing the edits in either order is safe -
# there's no stale-offset drift like there would be with Characters(start,len).

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(10)
$shape = $slide.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$target1 = "  plotHistogram(histogram, numNumbers);"
$replacement1 = "    plotHistogram(histogram, numNumbers);"
$found1 = $tr.Find($target1, 0)
if ($found1 -ne $null) {
    $found1.Text = $replacement1
} else {
    Write-Output "WARNING: could not find plotHistogram line"
}

$target2 = "  free(histogram); histogram = NULL;"
$replacement2 = "    free(histogram); histogram = NULL;"
$found2 = $tr.Find($target2, 0)
if ($found2 -ne $null) {
    $found2.Text = $replacement2
} else {
    Write-Output "WARNING: could not find free(histogram) line"
}
